$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.327.80'
$ws.Range('E2').Value = '  -2.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.735.52'
$ws.Range('E3').Value = '  -3.63%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.74'
$ws.Range('E5').Value = '  -3.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4237'
$ws.Range('E7').Value = '  -8.86%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3601'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.16'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07401'
$ws.Range('E10').Value = '  -3.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.111'
$ws.Range('E11').Value = '  -3.64%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.54'
$ws.Range('E13').Value = '  -4.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.065'
$ws.Range('E14').Value = '  -4.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.161'
$ws.Range('E15').Value = '  -3.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.729.36'
$ws.Range('E16').Value = '  -3.80%  '
$ws.Range('E17').Value = '  -3.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '86.98'
$ws.Range('E18').Value = '  +5.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06004'
$ws.Range('E19').Value = '  -10.77%  '
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.78'
$ws.Range('E21').Value = '  -3.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.068'
$ws.Range('E22').Value = '  -5.66%  '
$ws.Range('E23').Value = '  -3.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '27.363.06'
$ws.Range('E24').Value = '  -2.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.32'
$ws.Range('E25').Value = '  -5.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.384'
$ws.Range('E26').Value = '  -1.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.13'
$ws.Range('E27').Value = '  -3.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.364'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '149.04'
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.927.97'
$ws.Range('E30').Value = '  -3.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '126.14'
$ws.Range('E31').Value = '  -6.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.201'
$ws.Range('E32').Value = '  -5.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.631'
$ws.Range('E33').Value = '  -4.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09067'
$ws.Range('E34').Value = '  -5.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.626'
$ws.Range('E35').Value = '  -10.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.83'
$ws.Range('E36').Value = '  +4.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2132'
$ws.Range('E37').Value = '  -5.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.048'
$ws.Range('E38').Value = '  -4.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02247'
$ws.Range('E39').Value = '  -5.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06045'
$ws.Range('E40').Value = '  -5.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6363'
$ws.Range('E41').Value = '  -5.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.182'
$ws.Range('E42').Value = '  -4.28%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.928'
$ws.Range('E43').Value = '  -2.71%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9991'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.408'
$ws.Range('E45').Value = '  -7.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.57'
$ws.Range('E46').Value = '  -3.99%  '
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5817'
$ws.Range('E48').Value = '  -5.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.44'
$ws.Range('E49').Value = '  -4.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.951'
$ws.Range('E50').Value = '  -5.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06844'
$ws.Range('E51').Value = '  -4.17%  '
